$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "Ontology spreadsheet" value for row 579 (column K)
$ws.Range("K579").Value = "Intervention population"

# New rows of mapping data appended after row 579 (rows 580-587), columns A:K
$newRows = @(
    @("", "", "2", "", "tg1Ndrop", "Number of participants who dropped out of study arm (treatment group)", "GMHO:0000075", "number of participant drop-out from intervention", "Number of intervention participants who withdraw from or cannot complete an intervention.", "number of intervention participants", "Intervention outcomes and spillover effects"),
    @("", "", "2", "", "tg1Ndrop", "Number of participants who dropped out of study arm (treatment group)", "GMHO:0000123", "intervention arm", "A study arm that is subject to evaluation in an intervention evaluation study.", "study arm", "Intervention content and delivery"),
    @("", "", "2", "", "tg1Nrandom", "Number of intervention participants randomly allocated to an arm", "GMHO:0000154", "number of randomised intervention participants", "Number of intervention participants who were randomly allocated to study arms within an intervention.", "number of intervention participants", "Intervention population"),
    @("", "", "2", "", "tg1Nrandom", "Number of intervention participants randomly allocated to an arm", "GMHO:0000123", "intervention arm", "A study arm that is subject to evaluation in an intervention evaluation study.", "study arm", "Intervention content and delivery"),
    @("", "", "2", "", "tg2Ndrop", "Number of participants who dropped out of study arm (control group)", "GMHO:0000075", "number of participant drop-out from intervention", "Number of intervention participants who withdraw from or cannot complete an intervention.", "number of intervention participants", "Intervention outcomes and spillover effects"),
    @("", "", "2", "", "tg2Ndrop", "Number of participants who dropped out of study arm (control group)", "GMHO:0000122", "control arm", "A study arm designation as a comparator to some intervention arm.", "study arm", "Intervention content and delivery"),
    @("", "", "2", "", "tg2Nrandom", "Number of intervention participants randomly allocated to an arm", "GMHO:0000154", "number of randomised intervention participants", "Number of intervention participants who were randomly allocated to study arms within an intervention.", "number of intervention participants", "Intervention population"),
    @("", "", "2", "", "tg2Nrandom", "Number of intervention participants randomly allocated to an arm", "GMHO:0000122", "control arm", "A study arm designation as a comparator to some intervention arm.", "study arm", "Intervention content and delivery")
)

$startRow = 580
$endRow = $startRow + $newRows.Count - 1

# Column C ("LSR no.") holds numeric-looking text like "2" elsewhere in the sheet,
# so force it to text formatting before writing to preserve it as a string value.
$cFormatRange = "C" + $startRow + ":C" + $endRow
$ws.Range($cFormatRange).NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowData[$col - 1]
    }
}
